# Generate Report for Handoff
#
# "b.md" just finished its handoff cycle: a new handoff package was produced,
# so its status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", and the per-language sheets get a fresh handoff file
# name + handoff datetime recorded for the b.md row.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet: update the "b.md" row status for both languages ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# ---- zh-cn sheet: update status, new handoff file name, new handoff datetime ----
$wsZh = $wb.Worksheets.Item("zh-cn")

foreach ($link in $wsZh.Hyperlinks) {
    if ($link.Range.Address($false, $false) -eq "C3") {
        $link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-02-24 07:03:44"

# ---- de-de sheet: update status, new handoff file name, new handoff datetime ----
$wsDe = $wb.Worksheets.Item("de-de")

foreach ($link in $wsDe.Hyperlinks) {
    if ($link.Range.Address($false, $false) -eq "C3") {
        $link.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}

$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("D3").Value = "2016-02-24 07:03:58"
